$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 54, pushing existing rows 54..115 down to 55..116.
$ws.Rows("54:54").Insert()

# Populate the newly inserted row 54 with the new record.
$ws.Range("A54").Value = 11
$ws.Range("B54").Value = "Vega Monumental Concepción"
$ws.Range("C54").Value = "Bíobío"
$ws.Range("D54").Value = 44741
$ws.Range("D54").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E54").Value = 8
$ws.Range("F54").Value = 100112021
$ws.Range("G54").Value = "Ají"
$ws.Range("H54").Value = "Inferno"
$ws.Range("I54").Value = "Primera"
$ws.Range("J54").Value = 100
$ws.Range("K54").Value = 20000
$ws.Range("L54").Value = 22000
$ws.Range("M54").Value = 21000
$ws.Range("N54").Value = "$/caja 12 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 1750
$ws.Range("Q54").Value = 12
$ws.Range("R54").Value = "Hortaliza"
